# Add 2022-Q3 data
# 1) Insert a new summary row into "总计" sheet for the 2022-Q3 quarter.
# 2) Insert a brand-new "2022-Q3" worksheet (with the fund holdings detail)
#    right after "总计" and before "2022-Q2", pushing every later quarter
#    sheet one slot to the right.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row 2 for 2022-Q3
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
$summary.Range("A2:D2").ClearFormats()

# Re-apply the same style used by the other "A" column cells (bold/border)
$summary.Range("A4").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.29

# Column A is a 0-based row counter, independent of the quarter shown on
# each row; re-sequence it explicitly now that there are 8 data rows
# instead of 7 (rows 2..9 => 0..7).
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6
$summary.Range("A9").Value = 7

# ---------------------------------------------------------------------
# 2) New "2022-Q3" worksheet, cloned from the existing "2022-Q2" sheet so
#    it inherits identical column widths / styles, then overwritten with
#    the 2022-Q3 fund holdings.
#
# NOTE: worksheet variables captured *before* a Copy/Move that reshuffles
# tab order can end up referring to whichever sheet now sits at their old
# tab position rather than the sheet originally captured, so every sheet
# is (re)looked-up by name right before each reorder step below.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q2").Copy($null, $wb.Worksheets.Item("2022-Q2"))
$wb.Worksheets.Item(3).Name = "2022-Q3"
$wb.Worksheets.Item("2022-Q3").Move($wb.Worksheets.Item("2022-Q2"))

# Tab order is now settled; it is safe to keep a handle to the sheet from
# here on.
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

# The template (2022-Q2) only has 3 data rows; 2022-Q3 needs 4, so clone
# row 4's formatting down into a new row 5.
$q3Sheet.Rows.Item(4).Copy()
$q3Sheet.Rows.Item(5).Insert()
$q3Sheet.Range("A4").Copy()
$q3Sheet.Range("A5").PasteSpecial(-4122)

# Header row (unchanged labels, just making sure they are correct)
$q3Sheet.Range("B1").Value = "基金代码"
$q3Sheet.Range("C1").Value = "基金名称"
$q3Sheet.Range("D1").Value = "基金规模"
$q3Sheet.Range("E1").Value = "股票总仓位"
$q3Sheet.Range("F1").Value = "仓位占比"
$q3Sheet.Range("G1").Value = "持有市值(亿元)"
$q3Sheet.Range("H1").Value = "仓位排名"

# Row 2
$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("B2").Value = "'000586"
$q3Sheet.Range("C2").Value = "景顺长城中小创精选股票"
$q3Sheet.Range("D2").Value = "'2.21"
$q3Sheet.Range("E2").Value = "'93.50"
$q3Sheet.Range("F2").Value = "'7.94"
$q3Sheet.Range("G2").Value = "'0.1755"
$q3Sheet.Range("H2").Value = 3

# Row 3
$q3Sheet.Range("A3").Value = 1
$q3Sheet.Range("B3").Value = "'260115"
$q3Sheet.Range("C3").Value = "景顺长城中小盘混合"
$q3Sheet.Range("D3").Value = "'0.92"
$q3Sheet.Range("E3").Value = "'92.87"
$q3Sheet.Range("F3").Value = "'5.15"
$q3Sheet.Range("G3").Value = "'0.0474"
$q3Sheet.Range("H3").Value = 9

# Row 4
$q3Sheet.Range("A4").Value = 2
$q3Sheet.Range("B4").Value = "'010706"
$q3Sheet.Range("C4").Value = "景顺长城景骊成长混合"
$q3Sheet.Range("D4").Value = "'0.61"
$q3Sheet.Range("E4").Value = "'92.93"
$q3Sheet.Range("F4").Value = "'6.05"
$q3Sheet.Range("G4").Value = "'0.0369"
$q3Sheet.Range("H4").Value = 7

# Row 5
$q3Sheet.Range("A5").Value = 3
$q3Sheet.Range("B5").Value = "'000965"
$q3Sheet.Range("C5").Value = "汇丰晋信新动力混合"
$q3Sheet.Range("D5").Value = "'0.89"
$q3Sheet.Range("E5").Value = "'93.15"
$q3Sheet.Range("F5").Value = "'3.17"
$q3Sheet.Range("G5").Value = "'0.0282"
$q3Sheet.Range("H5").Value = 7

Write-Output "2022-Q3 sheet and summary row added."
